$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark. This removes both the
#    bookmarkStart that sat before "BAJO SEDACI..." and the
#    bookmarkEnd that sat right after " ESPASTICO )".
try {
    $d.Bookmarks("_GoBack").Delete()
} catch {
}

# 2. Locate the paragraph ending in "( COLON ESPASTICO )" - the block
#    of (originally 11) empty paragraphs that need trimming starts
#    right after it.
$anchor = $d.Content
$anchor.Find.Execute("( COLON ESPASTICO )", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$hostParaStart = $anchor.Paragraphs(1).Range.Start

$hostIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i++
    if ($p.Range.Start -eq $hostParaStart) {
        $hostIndex = $i
        break
    }
}

# Remove 5 of the empty paragraphs that sit right after that paragraph
# (there were 11 empty paragraphs; 5 get removed, leaving 6).
for ($k = 0; $k -lt 5; $k++) {
    $d.Paragraphs($hostIndex + 1).Range.Delete()
}

# 3. Re-add the "_GoBack" bookmark collapsed right after
#    "Dr. Gabriel Heriberto Gonzalez Asencio" (end of document
#    content, before the final paragraph mark). A collapsed bookmark
#    placed exactly at a paragraph-end boundary is not persisted by
#    this host, so insert a temporary placeholder character, wrap the
#    bookmark around it, then remove the placeholder again - the
#    bookmark collapses back in place and survives the save.
$endPos = $d.Content.End - 1
$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter([char]1)

$wrap = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $wrap)

$d.Range($endPos, $endPos + 1).Text = ""
